# Add the team's season record (Wins / Losses / Ties) as three new
# columns (AD, AE, AF) to the roster/statistics sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------
# Copy the formatting of the last existing header cell (AC1 - bold font,
# thin border, centered/top aligned) onto the three new header cells so
# they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-53) -------------------------------------------
# Every player row gets the same team record values.
$ws.Range("AD2:AD53").Value = 64
$ws.Range("AE2:AE53").Value = 98
$ws.Range("AF2:AF53").Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-53"
